# Update the Midway_B team transition-probability matrix on Sheet1.
# The simulation was re-run with more games, which shifted the
# empirical transition frequencies in several rows of the matrix.
# Each cell below is set to the recomputed probability value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2278481012658228
$ws.Range("C2").Value = 0.4683544303797468
$ws.Range("J2").Value = 0.08860759493670886
$ws.Range("P2").Value = 0.1645569620253164
$ws.Range("S2").Value = 0.05063291139240506
$ws.Range("J3").Value = 0.1351351351351351
$ws.Range("P3").Value = 0.6756756756756757
$ws.Range("S3").Value = 0.1891891891891892
$ws.Range("J4").Value = 0.1666666666666667
$ws.Range("P4").Value = 0.5833333333333334
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.06976744186046512
$ws.Range("D6").Value = 0.02325581395348837
$ws.Range("J6").Value = 0.2093023255813954
$ws.Range("Q6").Value = 0.2325581395348837
$ws.Range("R6").Value = 0.1162790697674419
$ws.Range("S6").Value = 0.3488372093023256
$ws.Range("B7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.2380952380952381
$ws.Range("Q7").Value = 0.2857142857142857
$ws.Range("R7").Value = 0.1428571428571428
$ws.Range("S7").Value = 0.2857142857142857
$ws.Range("B8").Value = 0.06024096385542169
$ws.Range("D8").Value = 0.02409638554216868
$ws.Range("F8").Value = 0.0963855421686747
$ws.Range("J8").Value = 0.144578313253012
$ws.Range("O8").Value = 0.01204819277108434
$ws.Range("Q8").Value = 0.1927710843373494
$ws.Range("R8").Value = 0.0963855421686747
$ws.Range("S8").Value = 0.3734939759036144
$ws.Range("B9").Value = 0.1
$ws.Range("F9").Value = 0.1
$ws.Range("J9").Value = 0.1
$ws.Range("O9").Value = 0.05
$ws.Range("Q9").Value = 0.175
$ws.Range("R9").Value = 0.075
$ws.Range("S9").Value = 0.4
$ws.Range("B10").Value = 0.1304347826086956
$ws.Range("D10").Value = 0.02608695652173913
$ws.Range("F10").Value = 0.06086956521739131
$ws.Range("J10").Value = 0.1217391304347826
$ws.Range("O10").Value = 0.002898550724637681
$ws.Range("Q10").Value = 0.2318840579710145
$ws.Range("R10").Value = 0.1217391304347826
$ws.Range("S10").Value = 0.3043478260869565
$ws.Range("G11").Value = 0.1724137931034483
$ws.Range("J11").Value = 0.1379310344827586
$ws.Range("K11").Value = 0.2758620689655172
$ws.Range("L11").Value = 0.3793103448275862
$ws.Range("S11").Value = 0.03448275862068965
$ws.Range("G12").Value = 0.8333333333333334
$ws.Range("L12").Value = 0.1666666666666667
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.2222222222222222
$ws.Range("S13").Value = 0.1111111111111111
$ws.Range("H15").Value = 0.1951219512195122
$ws.Range("I15").Value = 0.04878048780487805
$ws.Range("J15").Value = 0.5609756097560976
$ws.Range("S15").Value = 0.1951219512195122
$ws.Range("H16").Value = 0.1333333333333333
$ws.Range("I16").Value = 0.08888888888888889
$ws.Range("J16").Value = 0.6
$ws.Range("K16").Value = 0.04444444444444445
$ws.Range("O16").Value = 0.02222222222222222
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.01739130434782609
$ws.Range("H17").Value = 0.1043478260869565
$ws.Range("I17").Value = 0.05217391304347826
$ws.Range("J17").Value = 0.6434782608695652
$ws.Range("K17").Value = 0.03478260869565217
$ws.Range("M17").Value = 0.03478260869565217
$ws.Range("O17").Value = 0.04347826086956522
$ws.Range("S17").Value = 0.06956521739130435
$ws.Range("H18").Value = 0.1355932203389831
$ws.Range("I18").Value = 0.1016949152542373
$ws.Range("J18").Value = 0.5084745762711864
$ws.Range("M18").Value = 0.01694915254237288
$ws.Range("O18").Value = 0.1186440677966102
$ws.Range("S18").Value = 0.1186440677966102
$ws.Range("F19").Value = 0.01234567901234568
$ws.Range("H19").Value = 0.2098765432098765
$ws.Range("I19").Value = 0.08641975308641975
$ws.Range("J19").Value = 0.4197530864197531
$ws.Range("K19").Value = 0.05761316872427984
$ws.Range("M19").Value = 0.0205761316872428
$ws.Range("O19").Value = 0.07407407407407407
$ws.Range("S19").Value = 0.1193415637860082
